$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Array")
$ws.Name = "Arrays"
